# Weekly driver report update for 2025-04-28
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bad Drivers section
$ws.Range("C3").Value = 2844
$ws.Range("D3").Value = 79.59999999999999
$ws.Range("C4").Value = 2844

# Good Drivers section
$ws.Range("B12").Value = 11140
$ws.Range("B13").Value = 14487
